# atualizei dados bibi e add
# Insere 4 novos dias de faturamento de Agosto/2025 (dias 7-10) na planilha
# "faturamento_diario", empurrando os dados de Julho/Junho/Maio para baixo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# As linhas de dados comecam na linha 2; Agosto ja tem os dias 1-6 (linhas 2-7).
# Insere 4 linhas em branco logo depois (linhas 8-11) para os novos dias 7-10,
# deslocando Julho/Junho/Maio para baixo automaticamente.
$ws.Rows("8:11").Insert()

# Novos registros de Agosto/2025 (dias 7 a 10)
$novosDias = @(
    @{ Dia = 7;  Total = 28507.4;  Linha = 8 },
    @{ Dia = 8;  Total = 37782.45; Linha = 9 },
    @{ Dia = 9;  Total = 16187.29; Linha = 10 },
    @{ Dia = 10; Total = 2697.99;  Linha = 11 }
)

foreach ($registro in $novosDias) {
    $linha = $registro.Linha
    $ws.Range("A$linha").Value = $registro.Dia
    $ws.Range("B$linha").Value = $registro.Total
    $ws.Range("C$linha").Value = 8
    $ws.Range("D$linha").Value = 2025
    $ws.Range("E$linha").Value = "08/2025"
}
